$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update staff names to include the "Mrs." honorific (library/IQAC/staff data refresh)
$ws.Range("A12").Value = "Mrs. RHODAS DAISY D"
$ws.Range("A13").Value = "Mrs. GNANESHWARI R"

# Move the active selection to A13 (matches the saved selection state)
$ws.Range("A13").Select()
